$wb = $excel.ActiveWorkbook

# ---- Sheet "Overall" : row 2 updates ----
$wsOverall = $wb.Worksheets.Item("Overall")
$wsOverall.Range("B2").Value = 44
$wsOverall.Range("C2").Value = 5
$wsOverall.Range("D2").Value = 0.93731965832337605
$wsOverall.Range("E2").Value = 0.35
$wsOverall.Range("F2").Value = 1.0624861428840957
$wsOverall.Range("G2").Value = 27
$wsOverall.Range("H2").Value = 16
$wsOverall.Range("I2").Value = 43
$wsOverall.Range("J2").Value = 409
$wsOverall.Range("K2").Value = 18

# ---- Sheet "Zones" : rows 2-14 updates ----
$wsZones = $wb.Worksheets.Item("Zones")

# Row 2 (Zone 1)
$wsZones.Range("B2").Value = 5
$wsZones.Range("C2").Value = 0
$wsZones.Range("D2").Value = 0.92666666666666653
$wsZones.Range("E2").Value = 0.26666666666666661
$wsZones.Range("F2").Value = 1.0916666666666666

# Row 3 (Zone 2)
$wsZones.Range("B3").Value = 3
$wsZones.Range("C3").Value = 0
$wsZones.Range("D3").Value = 0.62083333333333335
$wsZones.Range("E3").Value = 0.32500000000000018
$wsZones.Range("F3").Value = 0.91666666666666663

# Row 4 (Zone 3)
$wsZones.Range("B4").Value = 9
$wsZones.Range("C4").Value = 2
$wsZones.Range("D4").Value = 0.83472222222222225
$wsZones.Range("E4").Value = 0.33333333333333304
$wsZones.Range("F4").Value = 0.93500000000000016

# Row 5 (Zone 4)
$wsZones.Range("B5").Value = 4
$wsZones.Range("C5").Value = 0
$wsZones.Range("D5").Value = 0.86041666666666672
$wsZones.Range("E5").Value = 0.15000000000000002
$wsZones.Range("F5").Value = 0.96190476190476193

# Row 6 (Zone 5) - no X_Wait (E) value
$wsZones.Range("B6").Value = 2
$wsZones.Range("C6").Value = 0
$wsZones.Range("D6").Value = 0.53333333333333321
$wsZones.Range("F6").Value = 0.53333333333333321

# Row 7 (Zone 6)
$wsZones.Range("B7").Value = 4
$wsZones.Range("C7").Value = 1
$wsZones.Range("D7").Value = 2.4537037037037042
$wsZones.Range("E7").Value = 0.21666666666666679
$wsZones.Range("F7").Value = 2.7333333333333334

# Row 8 (Zone 7) - no X_Wait (E) value
$wsZones.Range("B8").Value = 0
$wsZones.Range("C8").Value = 0
$wsZones.Range("D8").Value = 0.016666666666666607
$wsZones.Range("F8").Value = 0.016666666666666607

# Row 9 (Zone 8) - X_Wait (E) removed
$wsZones.Range("B9").Value = 2
$wsZones.Range("C9").Value = 0
$wsZones.Range("D9").Value = 0.49583333333333329
$wsZones.Range("E9").ClearContents()
$wsZones.Range("F9").Value = 0.49583333333333329

# Row 10 (Zone 9)
$wsZones.Range("B10").Value = 3
$wsZones.Range("C10").Value = 0
$wsZones.Range("D10").Value = 0.36666666666666686
$wsZones.Range("E10").Value = 0.41666666666666696
$wsZones.Range("F10").Value = 0.34166666666666679

# Row 11 (Zone 10)
$wsZones.Range("B11").Value = 4
$wsZones.Range("C11").Value = 0
$wsZones.Range("D11").Value = 0.5233333333333331
$wsZones.Range("E11").Value = 0.5583333333333329
$wsZones.Range("F11").Value = 0.49999999999999983

# Row 12 (Zone 11)
$wsZones.Range("B12").Value = 2
$wsZones.Range("C12").Value = 1
$wsZones.Range("D12").Value = 0.81904761904761891
$wsZones.Range("E12").Value = 0.18333333333333335
$wsZones.Range("F12").Value = 0.92499999999999982

# Row 13 (Zone 12) - no X_Wait (E) value
$wsZones.Range("B13").Value = 1
$wsZones.Range("C13").Value = 1
$wsZones.Range("D13").Value = 1.1638888888888888
$wsZones.Range("F13").Value = 1.1638888888888888

# Row 14 (Zone 13)
$wsZones.Range("B14").Value = 5
$wsZones.Range("C14").Value = 0
$wsZones.Range("D14").Value = 0.45833333333333354
$wsZones.Range("E14").Value = 0.4416666666666671
$wsZones.Range("F14").Value = 0.46666666666666679
